$p = $ppt.ActivePresentation

# --- 1. Slide 1 title text: "Zion" -> "Eden" ------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "The Innovation Zion Lifecycle Model") {
            $tr.Text = "The Innovation Eden Lifecycle Model"
        }
    }
}

# --- 2. Add a new 4th slide (Blank layout) --------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)

# Title textbox - duplicated (via copy/paste) from slide 1's title shape so
# that it keeps the exact accent2 / lumMod(60%) / lumOff(40%) font colour and
# the Magneto font formatting, then re-targeted to the new slide's size/text.
$titleSrc = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 3") {
        $titleSrc = $shp
        break
    }
}
$titleSrc.Copy()
$titlePasted = $newSlide.Shapes.Paste()
$titleBox = $titlePasted.Item(1)
$titleBox.Name = "TextBox 3"
$titleBox.Left = 0
$titleBox.Top = 0
$titleBox.Width = 960
$titleBox.Height = 60.58590551181102
$titleBox.TextFrame.TextRange.Text = "The Innovation Zion Optimal Process"

# "TBO!" textbox
$tboBox = $newSlide.Shapes.AddTextbox(1, 454.4276377952756, 255.45937007874016, 51.144645669291336, 29.081259842519685)
$tboBox.Name = "TextBox 4"
$tboBox.TextFrame.WordWrap = $false
$tboBox.TextFrame.AutoSize = 1
$tboBox.Fill.Visible = 0
$tboTr = $tboBox.TextFrame.TextRange
$tboTr.Text = "TBO!"
$tboTr.LanguageID = "en-GB"
$tboTr.ParagraphFormat.Alignment = 2
